# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet with the next data snapshot:
#  - bump the "Datos actualizados..." timestamp in A1
#  - update case/recovery/death figures for the countries whose numbers moved
#  - a handful of rows (129-138, 161-165) shifted position in the country
#    ranking, so both the country name (col A) and its stats (cols B-H)
#    change for those rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 7 de Abril de 2020 a las 00:22'
$ws.Range("B4").Value = 364013
$ws.Range("C4").Value = 27340
$ws.Range("D4").Value = 19532
$ws.Range("E4").Value = 333689
$ws.Range("F4").Value = 8878
$ws.Range("G4").Value = 1176
$ws.Range("H4").Value = 10792

$ws.Range("B7").Value = 102903
$ws.Range("C7").Value = 2780
$ws.Range("E7").Value = 72427
$ws.Range("G7").Value = 192
$ws.Range("H7").Value = 1776

$ws.Range("C9").Value = 0
$ws.Range("G9").Value = 0

$ws.Range("B16").Value = 16666
$ws.Range("C16").Value = 1154
$ws.Range("E16").Value = 12809
$ws.Range("G16").Value = 43
$ws.Range("H16").Value = 323

$ws.Range("E18").Value = 11367
$ws.Range("G18").Value = 76
$ws.Range("H18").Value = 562

$ws.Range("D43").Value = 118
$ws.Range("E43").Value = 2024

$ws.Range("B100").Value = 254
$ws.Range("C100").Value = 17
$ws.Range("E100").Value = 228

$ws.Range("A129").Value = 'Republica de Yibuti'
$ws.Range("B129").Value = 90
$ws.Range("C129").Value = 31
$ws.Range("D129").Value = 9
$ws.Range("E129").Value = 81
$ws.Range("F129").Value = 0

$ws.Range("A130").Value = 'Madagascar'
$ws.Range("B130").Value = 82
$ws.Range("C130").Value = 10
$ws.Range("D130").Value = 2
$ws.Range("E130").Value = 80
$ws.Range("F130").Value = 6
$ws.Range("H130").Value = 0

$ws.Range("A131").Value = 'Monaco'
$ws.Range("C131").Value = 4
$ws.Range("D131").Value = 4
$ws.Range("E131").Value = 72
$ws.Range("F131").Value = 4

$ws.Range("A132").Value = 'Liechtenstein'
$ws.Range("B132").Value = 77
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 55
$ws.Range("E132").Value = 21
$ws.Range("H132").Value = 1

$ws.Range("A133").Value = 'Aruba'
$ws.Range("B133").Value = 71
$ws.Range("C133").Value = 7
$ws.Range("D133").Value = 2
$ws.Range("E133").Value = 69
$ws.Range("F133").Value = 0
$ws.Range("H133").Value = 0

$ws.Range("A134").Value = 'Guatemala'
$ws.Range("B134").Value = 70
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 15
$ws.Range("E134").Value = 52
$ws.Range("F134").Value = 3
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 3

$ws.Range("A135").Value = 'El Salvador'
$ws.Range("B135").Value = 69
$ws.Range("C135").Value = 7
$ws.Range("D135").Value = 5
$ws.Range("E135").Value = 60
$ws.Range("F135").Value = 4
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 4

$ws.Range("A136").Value = 'Guayana Francesa'
$ws.Range("B136").Value = 68
$ws.Range("D136").Value = 27
$ws.Range("E136").Value = 41
$ws.Range("F136").Value = 1

$ws.Range("A137").Value = 'Barbados'
$ws.Range("B137").Value = 60
$ws.Range("C137").Value = 4
$ws.Range("D137").Value = 6
$ws.Range("E137").Value = 52
$ws.Range("F137").Value = 4
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 2

$ws.Range("A138").Value = 'Jamaica'
$ws.Range("B138").Value = 58
$ws.Range("D138").Value = 8
$ws.Range("E138").Value = 47
$ws.Range("F138").Value = 0
$ws.Range("H138").Value = 3

$ws.Range("B154").Value = 31
$ws.Range("C154").Value = 7
$ws.Range("D154").Value = 8
$ws.Range("E154").Value = 19
$ws.Range("F154").Value = 8

$ws.Range("A161").Value = 'Libia'
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 1
$ws.Range("E161").Value = 17
$ws.Range("H161").Value = 1

$ws.Range("A162").Value = 'Siria'
$ws.Range("D162").Value = 2
$ws.Range("E162").Value = 15
$ws.Range("H162").Value = 2

$ws.Range("A163").Value = 'Maldivas'
$ws.Range("B163").Value = 19
$ws.Range("D163").Value = 13
$ws.Range("E163").Value = 6

$ws.Range("A164").Value = 'Guinea-Bisau'
$ws.Range("D164").Value = 0
$ws.Range("E164").Value = 18

$ws.Range("A165").Value = 'Nueva Caledonia'
$ws.Range("E165").Value = 17
$ws.Range("H165").Value = 0
